# Auto-generated edit script: applies updated market-price/profit values
# to the Leve profit tables across multiple sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 23
$ws.Range("I8").Value = 23
$ws.Range("K8").Value = 69
$ws.Range("M8").Value = 70
$ws.Range("H62").Value = 4066.6667
$ws.Range("I62").Value = 3701.3333
$ws.Range("J62").Value = 4432
$ws.Range("K62").Value = 3701.3333
$ws.Range("L62").Value = 4432
$ws.Range("M62").Value = -3077.3333
$ws.Range("N62").Value = -5680
$ws.Range("H65").Value = 4066.6667
$ws.Range("I65").Value = 3701.3333
$ws.Range("J65").Value = 4432
$ws.Range("K65").Value = 18506.6665
$ws.Range("L65").Value = 22160
$ws.Range("M65").Value = -15386.6665
$ws.Range("N65").Value = -28400
$ws.Range("H70").Value = 4999.5
$ws.Range("J70").Value = 4999.5
$ws.Range("L70").Value = 14998.5
$ws.Range("N70").Value = -15538.5
$ws.Range("H73").Value = 4999.5
$ws.Range("J73").Value = 4999.5
$ws.Range("L73").Value = 14998.5
$ws.Range("N73").Value = -16870.5
$ws.Range("H87").Value = 32500
$ws.Range("J87").Value = 32500
$ws.Range("L87").Value = 32500
$ws.Range("N87").Value = -34996
$ws.Range("H90").Value = 32500
$ws.Range("J90").Value = 32500
$ws.Range("L90").Value = 97500
$ws.Range("N90").Value = -109980
$ws.Range("H98").Value = 2904.2
$ws.Range("I98").Value = 2796.5833
$ws.Range("K98").Value = 2796.5833
$ws.Range("M98").Value = -1298.5833
$ws.Range("H113").Value = 4831.3335
$ws.Range("I113").Value = 4831.3335
$ws.Range("K113").Value = 4831.3335
$ws.Range("M113").Value = -1577.3335
$ws.Range("H122").Value = 2904.2
$ws.Range("I122").Value = 2796.5833
$ws.Range("K122").Value = 8389.749899999999
$ws.Range("M122").Value = -5939.749899999999
$ws.Range("H132").Value = 1895.3636
$ws.Range("J132").Value = 2437.5715
$ws.Range("L132").Value = 7312.7145
$ws.Range("N132").Value = -12372.7145
$ws.Range("H137").Value = 2072.7273
$ws.Range("I137").Value = 1790.4445
$ws.Range("J137").Value = 3343
$ws.Range("K137").Value = 5371.333500000001
$ws.Range("L137").Value = 10029
$ws.Range("M137").Value = -2821.333500000001
$ws.Range("N137").Value = -15129

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 1000
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1340
$ws.Range("H94").Value = 2131.1
$ws.Range("I94").Value = 2365.75
$ws.Range("K94").Value = 2365.75
$ws.Range("M94").Value = -1914.75
$ws.Range("H107").Value = 1715.4546
$ws.Range("I107").Value = 608.875
$ws.Range("K107").Value = 608.875
$ws.Range("M107").Value = 1311.125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 208
$ws.Range("I2").Value = 139.44444
$ws.Range("J2").Value = 362.25
$ws.Range("K2").Value = 139.44444
$ws.Range("L2").Value = 362.25
$ws.Range("M2").Value = -26.44443999999999
$ws.Range("N2").Value = -588.25
$ws.Range("H17").Value = 200
$ws.Range("I17").Value = 200
$ws.Range("K17").Value = 200
$ws.Range("M17").Value = -26
$ws.Range("H22").Value = 229.1875
$ws.Range("I22").Value = 211.3077
$ws.Range("J22").Value = 306.66666
$ws.Range("K22").Value = 211.3077
$ws.Range("L22").Value = 306.66666
$ws.Range("M22").Value = 138.6923
$ws.Range("N22").Value = -1006.66666
$ws.Range("H41").Value = 16250
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 20000
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 20000
$ws.Range("M41").Value = -4572
$ws.Range("N41").Value = -20856
$ws.Range("H86").Value = 9250
$ws.Range("I86").Value = 9250
$ws.Range("K86").Value = 9250
$ws.Range("M86").Value = -8127
$ws.Range("H88").Value = 14570.857
$ws.Range("J88").Value = 14570.857
$ws.Range("L88").Value = 14570.857
$ws.Range("N88").Value = -15382.857
$ws.Range("H89").Value = 9250
$ws.Range("I89").Value = 9250
$ws.Range("K89").Value = 46250
$ws.Range("M89").Value = -40634
$ws.Range("H91").Value = 14570.857
$ws.Range("J91").Value = 14570.857
$ws.Range("L91").Value = 14570.857
$ws.Range("N91").Value = -17378.857
$ws.Range("H105").Value = 1114.8334
$ws.Range("I105").Value = 1201
$ws.Range("K105").Value = 1201
$ws.Range("M105").Value = 546
$ws.Range("H107").Value = 555.8
$ws.Range("I107").Value = 750
$ws.Range("J107").Value = 507.25
$ws.Range("K107").Value = 750
$ws.Range("L107").Value = 507.25
$ws.Range("M107").Value = 1170
$ws.Range("N107").Value = -4347.25
$ws.Range("H134").Value = 1950.1818
$ws.Range("I134").Value = 1103.2354
$ws.Range("J134").Value = 4829.8
$ws.Range("K134").Value = 3309.7062
$ws.Range("L134").Value = 14489.4
$ws.Range("M134").Value = -774.7062000000001
$ws.Range("N134").Value = -19559.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 101
$ws.Range("I19").Value = 100
$ws.Range("K19").Value = 300
$ws.Range("M19").Value = -126
$ws.Range("H109").Value = 1914
$ws.Range("I109").Value = 1931.5
$ws.Range("K109").Value = 5794.5
$ws.Range("M109").Value = -4754.5
$ws.Range("H132").Value = 1749.5
$ws.Range("J132").Value = 1500
$ws.Range("L132").Value = 13500
$ws.Range("N132").Value = -18560

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3527.6538
$ws.Range("I80").Value = 1804.6666
$ws.Range("J80").Value = 4044.55
$ws.Range("K80").Value = 1804.6666
$ws.Range("L80").Value = 4044.55
$ws.Range("M80").Value = -806.6666
$ws.Range("N80").Value = -6040.55
$ws.Range("H83").Value = 3527.6538
$ws.Range("I83").Value = 1804.6666
$ws.Range("J83").Value = 4044.55
$ws.Range("K83").Value = 9023.333000000001
$ws.Range("L83").Value = 20222.75
$ws.Range("M83").Value = -4031.333000000001
$ws.Range("N83").Value = -30206.75
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1888.7778
$ws.Range("I22").Value = 1499.5
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1499.5
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -1204.5
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 1888.7778
$ws.Range("I27").Value = 1499.5
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 1499.5
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -1392.5
$ws.Range("N27").Value = -2214
$ws.Range("H46").Value = 3461.3845
$ws.Range("I46").Value = 2833.3333
$ws.Range("J46").Value = 3649.8
$ws.Range("K46").Value = 2833.3333
$ws.Range("L46").Value = 3649.8
$ws.Range("M46").Value = -2645.3333
$ws.Range("N46").Value = -4025.8
$ws.Range("H55").Value = 4025.8
$ws.Range("I55").Value = 2575
$ws.Range("J55").Value = 4993
$ws.Range("K55").Value = 2575
$ws.Range("L55").Value = 4993
$ws.Range("M55").Value = -2402
$ws.Range("N55").Value = -5339
$ws.Range("H68").Value = 3883.3333
$ws.Range("I68").Value = 3883.3333
$ws.Range("K68").Value = 3883.3333
$ws.Range("M68").Value = -3134.3333
$ws.Range("H71").Value = 3883.3333
$ws.Range("I71").Value = 3883.3333
$ws.Range("K71").Value = 19416.6665
$ws.Range("M71").Value = -15672.6665
$ws.Range("H82").Value = 2014.5
$ws.Range("I82").Value = 2037
$ws.Range("J82").Value = 1999.5
$ws.Range("K82").Value = 2037
$ws.Range("L82").Value = 1999.5
$ws.Range("M82").Value = -1676
$ws.Range("N82").Value = -2721.5
$ws.Range("H85").Value = 2014.5
$ws.Range("I85").Value = 2037
$ws.Range("J85").Value = 1999.5
$ws.Range("K85").Value = 2037
$ws.Range("L85").Value = 1999.5
$ws.Range("M85").Value = -789
$ws.Range("N85").Value = -4495.5
$ws.Range("H122").Value = 3215.5
$ws.Range("I122").Value = 3317.5454
$ws.Range("J122").Value = 2991
$ws.Range("K122").Value = 9952.636200000001
$ws.Range("L122").Value = 8973
$ws.Range("M122").Value = -7502.636200000001
$ws.Range("N122").Value = -13873
$ws.Range("H132").Value = 9439.143
$ws.Range("I132").Value = 9745.362999999999
$ws.Range("J132").Value = 8316.333000000001
$ws.Range("K132").Value = 29236.089
$ws.Range("L132").Value = 24948.999
$ws.Range("M132").Value = -26706.089
$ws.Range("N132").Value = -30008.999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("H132").Value = 1139
$ws.Range("J132").Value = 1350
$ws.Range("L132").Value = 4050
$ws.Range("N132").Value = -9110
